# Update forecast values ("预测值", column C) on Sheet1 for rows 2-24
# per the latest model run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = -52.8157
    3  = -29.3906
    4  = -60.8884
    5  = -76.05629999999999
    6  = -70.2509
    7  = -147.3776
    8  = -162.2014
    9  = -146.154
    10 = -80.1519
    11 = -30.3341
    12 = -78.00060000000001
    13 = -78.9799
    14 = -78.5346
    15 = -29.1136
    16 = -50.8619
    17 = -94.9991
    18 = -46.1783
    19 = -1.7846
    20 = 62.7789
    21 = 56.6976
    22 = 17.7444
    23 = 58.0395
    24 = 113.8637
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}
